# Turn the plain-text GitHub URL on the "GitHub Link" slide into a live
# hyperlink, splitting it the way PowerPoint does when you select the text
# and run Insert > Link (the trailing ".git" ends up as its own run), and
# leave an extra empty paragraph behind the link line.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$url = "https://github.com/sanketrshinde19/-Steganography.git"

# Make sure we're starting from the expected plain-text run before editing.
if ($tr.Text -eq $url) {

    # Press Enter at the end of the line first (matches the extra blank
    # paragraph seen after the link in the final slide), then hyperlink the
    # URL text itself in two pieces - "...Steganography" and ".git" - both
    # pointing at the same address.
    [void]$tr.InsertAfter("`r")

    $urlPart1Len = $url.Length - 4  # everything up to ".git"
    $part1 = $tr.Characters(1, $urlPart1Len)
    $part1.ActionSettings.Item(1).Hyperlink.Address = $url

    $part2 = $tr.Characters($urlPart1Len + 1, 4)
    $part2.ActionSettings.Item(1).Hyperlink.Address = $url
}
